$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.622.23"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.927.26"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  +0.42%  "
$ws.Range("D5").Value = "326.56"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").Value = "1.011"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").Value = "0.4827"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "0.4066"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "0.08236"
$ws.Range("E9").Value = "  +1.24%  "
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("D11").Value = "23.91"
$ws.Range("E11").Value = "  +2.09%  "
$ws.Range("D12").Value = "1.936.16"
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("D13").Value = "6.113"
$ws.Range("E13").Value = "  +2.15%  "
$ws.Range("D14").Value = "7.277"
$ws.Range("E14").Value = "  +2.23%  "
$ws.Range("D15").Value = "91.96"
$ws.Range("E15").Value = "  +2.04%  "
$ws.Range("D16").Value = "0.06886"
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("D17").Value = "1.013"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "17.68"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "1.011"
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("D21").Value = "29.626.17"
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("D22").Value = "5.692"
$ws.Range("E22").Value = "  +1.19%  "
$ws.Range("D24").Value = "2.188"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").Value = "2.162.03"
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("D26").Value = "156.38"
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("D27").Value = "6.457"
$ws.Range("E27").Value = "  +1.11%  "
$ws.Range("D28").Value = "20.08"
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("D29").Value = "2.099"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "120.77"
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("D31").Value = "1.016"
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("D32").Value = "0.09646"
$ws.Range("E32").Value = "  +1.19%  "
$ws.Range("Z1").Formula = "=`"5.630`""
$ws.Range("Z1").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E33").Value = "  +2.22%  "
$ws.Range("D34").Value = "3.578"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("D35").Value = "1.383"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").Value = "0.06407"
$ws.Range("E36").Value = "  +5.21%  "
$ws.Range("D37").Value = "0.02297"
$ws.Range("E37").Value = "  +1.35%  "
$ws.Range("E38").Value = "  +1.10%  "
$ws.Range("Z1").Formula = "=`"0.5970`""
$ws.Range("Z1").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("E39").Value = "  +0.54%  "
$ws.Range("D40").Value = "10.76"
$ws.Range("E40").Value = "  +0.87%  "
$ws.Range("D41").Value = "7.893"
$ws.Range("E41").Value = "  -0.65%  "
$ws.Range("D42").Value = "0.1855"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").Value = "2.454"
$ws.Range("E43").Value = "  +1.94%  "
$ws.Range("D44").Value = "1.285"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").Value = "12.42"
$ws.Range("E45").Value = "  -0.97%  "
$ws.Range("D46").Value = "0.07552"
$ws.Range("D47").Value = "0.5576"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("D48").Value = "1.964"
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("D49").Value = "119.62"
$ws.Range("E49").Value = "  +2.92%  "
$ws.Range("D50").Value = "2.443"
$ws.Range("E50").Value = "  +3.72%  "
$ws.Range("D51").Value = "72.28"
$ws.Range("E51").Value = "  -0.18%  "
